$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell text (Drools decision-table rewrite) ---

# Row 1: Import package cell (C1): com.myspace.eotworkflow (unchanged text,
# but the shared-string table is being reshuffled so just re-assert values).
$ws.Range("C1").Value = 'com.myspace.eotworkflow'

# Row 2: Import statement + wildcard import.
$ws.Range("B2").Value = 'Import'
$ws.Range("C2").Value = 'com.myspace.eotworkflow.*'

# Row 3 (new row, previously empty): Declare / dialect "mvel"; header that
# was pasted in above the rule table, using the "Normal 2" / Arial 10pt
# style that comes along with it.
$ws.Range("B3").Value = 'Declare'
$ws.Range("C3").Value = 'dialect "mvel";'
$ws.Range("B3:C3").Font.Name = "Arial"
$ws.Range("B3:C3").Font.Size = 10

# Row 4: RuleTable header label stays the same text.
$ws.Range("B4").Value = 'RuleTable AutoDocGen'

# Row 6: object declaration switches from the old "$document" MVEL binding
# syntax to the new bare "document" binding (works with the Declare/dialect
# header added above).
$ws.Range("B6").Value = 'document: Document'

# Row 7: column headers for the object's fields.
$ws.Range("B7").Value = 'relationshipname'
$ws.Range("C7").Value = 'docname'
$ws.Range("D7").Value = 'desc'
$ws.Range("F7").Value = 'document.setAutogen($param);'

# Row 8: descriptive sub-headers.
$ws.Range("B8").Value = 'relationship name'
$ws.Range("C8").Value = 'doc name'
$ws.Range("D8").Value = 'desc'

# Row 9: first data row.
$ws.Range("B9").Value = 'Buyout'
$ws.Range("C9").Value = 'OffertoPurchase'
$ws.Range("D9").Value = 'EOT'
$ws.Range("E9").Value = 'Auto Doc Gen'

# Row 10: second data row.
$ws.Range("B10").Value = 'Buyout'
$ws.Range("C10").Value = 'Bill of Sale'
$ws.Range("D10").Value = 'EOT'
$ws.Range("E10").Value = 'Auto Doc Gen'
